$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (Target cluster = ECs)
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.05619066666666667
$ws.Range("N2").Value = 0.168572
$ws.Range("O2").Value = 0.3931387525216601
$ws.Range("P2").Value = 0.39313875252166
$ws.Range("Q2").Value = 0.1037026099457778
$ws.Range("R2").Value = 0.9333234895120001
$ws.Range("S2").Value = 0.3931387525216601
$ws.Range("T2").Value = 0.39313875252166

# Update row 3 values (Target cluster = FAPs)
$ws.Range("O3").Value = 0.60686124747834
$ws.Range("P3").Value = 0.60686124747834
$ws.Range("S3").Value = 0.60686124747834
$ws.Range("T3").Value = 0.60686124747834

# Delete row 4 entirely (Target cluster = MuSCs)
$ws.Range("A4:T4").EntireRow.Delete()
